$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5: terminology change (Stationaer/Erledigt/Ambulant -> IMP/finished/AMB) ---
$ws.Range("J2").Value = "IMP"
$ws.Range("K2").Value = "finished"
$ws.Range("J3").Value = "AMB"
$ws.Range("K3").Value = "finished"
$ws.Range("J4").Value = "IMP"
$ws.Range("K4").Value = "finished"
$ws.Range("J5").Value = "AMB"
$ws.Range("K5").Value = "finished"

# --- Rows 6-30: add J:L (IMP/finished/HD), M (date), and N:O where applicable ---
$ws.Range("J6").Value = "IMP"
$ws.Range("K6").Value = "finished"
$ws.Range("L6").Value = "HD"
$ws.Range("M6").Value = 43561
$ws.Range("M6").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N6").Value = "DE"
$ws.Range("O6").Value = 4503

$ws.Range("J7").Value = "IMP"
$ws.Range("K7").Value = "finished"
$ws.Range("L7").Value = "HD"
$ws.Range("M7").Value = 43562
$ws.Range("M7").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("J8").Value = "IMP"
$ws.Range("K8").Value = "finished"
$ws.Range("L8").Value = "HD"
$ws.Range("M8").Value = 43563
$ws.Range("M8").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("J9").Value = "IMP"
$ws.Range("K9").Value = "finished"
$ws.Range("L9").Value = "HD"
$ws.Range("M9").Value = 43564
$ws.Range("M9").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("J10").Value = "IMP"
$ws.Range("K10").Value = "finished"
$ws.Range("L10").Value = "HD"
$ws.Range("M10").Value = 43565
$ws.Range("M10").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N10").Value = "DE"
$ws.Range("O10").Value = 5303

$ws.Range("J11").Value = "IMP"
$ws.Range("K11").Value = "finished"
$ws.Range("L11").Value = "HD"
$ws.Range("M11").Value = 43566
$ws.Range("M11").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N11").Value = "DE"
$ws.Range("O11").Value = 5503

$ws.Range("J12").Value = "IMP"
$ws.Range("K12").Value = "finished"
$ws.Range("L12").Value = "HD"
$ws.Range("M12").Value = 43567
$ws.Range("M12").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N12").Value = "DE"
$ws.Range("O12").Value = 5703

$ws.Range("J13").Value = "IMP"
$ws.Range("K13").Value = "finished"
$ws.Range("L13").Value = "HD"
$ws.Range("M13").Value = 43568
$ws.Range("M13").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N13").Value = "DE"
$ws.Range("O13").Value = 5903

$ws.Range("J14").Value = "IMP"
$ws.Range("K14").Value = "finished"
$ws.Range("L14").Value = "HD"
$ws.Range("M14").Value = 43569
$ws.Range("M14").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N14").Value = "DE"
$ws.Range("O14").Value = 6103

$ws.Range("J15").Value = "IMP"
$ws.Range("K15").Value = "finished"
$ws.Range("L15").Value = "HD"
$ws.Range("M15").Value = 43570
$ws.Range("M15").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N15").Value = "DE"
$ws.Range("O15").Value = 6303

$ws.Range("J16").Value = "IMP"
$ws.Range("K16").Value = "finished"
$ws.Range("L16").Value = "HD"
$ws.Range("M16").Value = 43571
$ws.Range("M16").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N16").Value = "DE"
$ws.Range("O16").Value = 6503

$ws.Range("J17").Value = "IMP"
$ws.Range("K17").Value = "finished"
$ws.Range("L17").Value = "HD"
$ws.Range("M17").Value = 43572
$ws.Range("M17").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N17").Value = "DE"
$ws.Range("O17").Value = 6703

$ws.Range("J18").Value = "IMP"
$ws.Range("K18").Value = "finished"
$ws.Range("L18").Value = "HD"
$ws.Range("M18").Value = 43573
$ws.Range("M18").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N18").Value = "DE"
$ws.Range("O18").Value = 6903

$ws.Range("J19").Value = "IMP"
$ws.Range("K19").Value = "finished"
$ws.Range("L19").Value = "HD"
$ws.Range("M19").Value = 43574
$ws.Range("M19").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N19").Value = "DE"
$ws.Range("O19").Value = 7103

$ws.Range("J20").Value = "IMP"
$ws.Range("K20").Value = "finished"
$ws.Range("L20").Value = "HD"
$ws.Range("M20").Value = 43575
$ws.Range("M20").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N20").Value = "DE"
$ws.Range("O20").Value = 7303

$ws.Range("J21").Value = "IMP"
$ws.Range("K21").Value = "finished"
$ws.Range("L21").Value = "HD"
$ws.Range("M21").Value = 43576
$ws.Range("M21").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N21").Value = "DE"
$ws.Range("O21").Value = 7503

$ws.Range("J22").Value = "IMP"
$ws.Range("K22").Value = "finished"
$ws.Range("L22").Value = "HD"
$ws.Range("M22").Value = 43577
$ws.Range("M22").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N22").Value = "DE"
$ws.Range("O22").Value = 7703

$ws.Range("J23").Value = "IMP"
$ws.Range("K23").Value = "finished"
$ws.Range("L23").Value = "HD"
$ws.Range("M23").Value = 43578
$ws.Range("M23").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N23").Value = "DE"
$ws.Range("O23").Value = 7903

$ws.Range("J24").Value = "IMP"
$ws.Range("K24").Value = "finished"
$ws.Range("L24").Value = "HD"
$ws.Range("M24").Value = 43579
$ws.Range("M24").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N24").Value = "DE"
$ws.Range("O24").Value = 8103

$ws.Range("J25").Value = "IMP"
$ws.Range("K25").Value = "finished"
$ws.Range("L25").Value = "HD"
$ws.Range("M25").Value = 43580
$ws.Range("M25").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N25").Value = "DE"
$ws.Range("O25").Value = 8303

$ws.Range("J26").Value = "IMP"
$ws.Range("K26").Value = "finished"
$ws.Range("L26").Value = "HD"
$ws.Range("M26").Value = 43581
$ws.Range("M26").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N26").Value = "DE"
$ws.Range("O26").Value = 8503

$ws.Range("J27").Value = "IMP"
$ws.Range("K27").Value = "finished"
$ws.Range("L27").Value = "HD"
$ws.Range("M27").Value = 43582
$ws.Range("M27").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N27").Value = "DE"
$ws.Range("O27").Value = 8703

$ws.Range("J28").Value = "IMP"
$ws.Range("K28").Value = "finished"
$ws.Range("L28").Value = "HD"
$ws.Range("M28").Value = 43583
$ws.Range("M28").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N28").Value = "DE"
$ws.Range("O28").Value = 8903

$ws.Range("J29").Value = "IMP"
$ws.Range("K29").Value = "finished"
$ws.Range("L29").Value = "HD"
$ws.Range("M29").Value = 43584
$ws.Range("M29").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N29").Value = "DE"
$ws.Range("O29").Value = 9103

$ws.Range("J30").Value = "IMP"
$ws.Range("K30").Value = "finished"
$ws.Range("L30").Value = "HD"
$ws.Range("M30").Value = 43585
$ws.Range("M30").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N30").Value = "DE"
$ws.Range("O30").Value = 9303

# --- New rows 31 and 32 (duplicates of row 26 data, with J:O additions) ---
$ws.Range("A31").Value = "259294944-TestHaus"
$ws.Range("B31").Value = "P_20085770"
$ws.Range("C31").Value = "F_101664"
$ws.Range("D31").Value = 1990
$ws.Range("E31").Value = "female"
$ws.Range("F31").Value = 55155
$ws.Range("F31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("G31").Value = 43923
$ws.Range("G31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("H31").Value = "E75.0"
$ws.Range("I31").Value = 796
$ws.Range("J31").Value = "IMP"
$ws.Range("K31").Value = "finished"
$ws.Range("L31").Value = "HD"
$ws.Range("M31").Value = 43581
$ws.Range("M31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N31").Value = "DE"
$ws.Range("O31").Value = 8503

$ws.Range("A32").Value = "259294944-TestHaus"
$ws.Range("B32").Value = "P_20085770"
$ws.Range("C32").Value = "F_101664"
$ws.Range("D32").Value = 1990
$ws.Range("E32").Value = "female"
$ws.Range("F32").Value = 55155
$ws.Range("F32").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("G32").Value = 43923
$ws.Range("G32").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("H32").Value = "E75.0"
$ws.Range("I32").Value = 796
$ws.Range("J32").Value = "IMP"
$ws.Range("K32").Value = "finished"
$ws.Range("L32").Value = "HD"
$ws.Range("M32").Value = 43581
$ws.Range("M32").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("N32").Value = "DE"
$ws.Range("O32").Value = 8503

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 110
$ws.Range("A34").Select() | Out-Null

# --- Column G (index 7) width: best-effort closest achievable to 16.6 ---
$ws.Columns.Item(7).ColumnWidth = 15.75
